$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 50001468
$ws.Range("I17").Value = 750
$ws.Range("J17").Value = 62501650
$ws.Range("K17").Value = 2250
$ws.Range("L17").Value = 187504950
$ws.Range("M17").Value = -2082
$ws.Range("N17").Value = -187505286

$ws.Range("H33").Value = 1096.1111
$ws.Range("I33").Value = 1074.2858
$ws.Range("J33").Value = 1172.5
$ws.Range("K33").Value = 1074.2858
$ws.Range("L33").Value = 1172.5
$ws.Range("M33").Value = -845.2858000000001
$ws.Range("N33").Value = -1630.5

$ws.Range("H112").Value = 1629.8368
$ws.Range("J112").Value = 1633.0625
$ws.Range("L112").Value = 4899.1875
$ws.Range("N112").Value = -7115.1875

$ws.Range("H116").Value = 3263
$ws.Range("I116").Value = 3263
$ws.Range("K116").Value = 3263
$ws.Range("M116").Value = 179

$ws.Range("H124").Value = 59980
$ws.Range("J124").Value = 59980
$ws.Range("L124").Value = 59980
$ws.Range("N124").Value = -69800

$ws.Range("H138").Value = 2563.3447
$ws.Range("I138").Value = 984.7646999999999
$ws.Range("K138").Value = 2954.2941
$ws.Range("M138").Value = 2185.7059

$ws.Range("H141").Value = 4732.154
$ws.Range("I141").Value = 4719.8335
$ws.Range("K141").Value = 14159.5005
$ws.Range("M141").Value = -8979.500499999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3583.2942
$ws.Range("I74").Value = 3453
$ws.Range("K74").Value = 3453
$ws.Range("M74").Value = -2579

$ws.Range("H77").Value = 3583.2942
$ws.Range("I77").Value = 3453
$ws.Range("K77").Value = 17265
$ws.Range("M77").Value = -12897

$ws.Range("H110").Value = 1817.9546
$ws.Range("I110").Value = 2011.0555
$ws.Range("K110").Value = 2011.0555
$ws.Range("M110").Value = 33.94450000000006

$ws.Range("H132").Value = 41671084
$ws.Range("I132").Value = 76926984
$ws.Range("K132").Value = 230780952
$ws.Range("M132").Value = -230778422

$ws.Range("H133").Value = 68600.91

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 49999.5
$ws.Range("J126").Value = 49999.5
$ws.Range("L126").Value = 49999.5
$ws.Range("N126").Value = -59879.5

$ws.Range("H132").Value = 70153.766
$ws.Range("J132").Value = 70153.766
$ws.Range("L132").Value = 70153.766
$ws.Range("N132").Value = -80273.766

$ws.Range("H134").Value = 3618.9033
$ws.Range("I134").Value = 2041.5238
$ws.Range("K134").Value = 6124.5714
$ws.Range("M134").Value = -3589.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 329
$ws.Range("I22").Value = 329
$ws.Range("K22").Value = 329
$ws.Range("M22").Value = 21

$ws.Range("H31").Value = 3865.6128
$ws.Range("I31").Value = 1818
$ws.Range("K31").Value = 1818
$ws.Range("M31").Value = -1523

$ws.Range("H34").Value = 3865.6128
$ws.Range("I34").Value = 1818
$ws.Range("K34").Value = 1818
$ws.Range("M34").Value = -1616

$ws.Range("H132").Value = 2396.25
$ws.Range("I132").Value = 1922.1578
$ws.Range("K132").Value = 5766.4734
$ws.Range("M132").Value = -3236.4734

$ws.Range("H134").Value = 3523.8135
$ws.Range("I134").Value = 2463.2
$ws.Range("K134").Value = 7389.599999999999
$ws.Range("M134").Value = -4854.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5171.5
$ws.Range("I3").Value = 5343
$ws.Range("K3").Value = 16029
$ws.Range("M3").Value = -15917

$ws.Range("H14").Value = 699.75
$ws.Range("I14").Value = 699.75
$ws.Range("K14").Value = 2099.25
$ws.Range("M14").Value = -1926.25

$ws.Range("H51").Value = 1550
$ws.Range("I51").Value = 1550
$ws.Range("K51").Value = 4650
$ws.Range("M51").Value = -4190

$ws.Range("H56").Value = 7864.9287
$ws.Range("I56").Value = 7864.9287
$ws.Range("K56").Value = 7864.9287
$ws.Range("M56").Value = -7334.9287

$ws.Range("H60").Value = 1201.3077
$ws.Range("I60").Value = 171.8
$ws.Range("J60").Value = 1844.75
$ws.Range("K60").Value = 515.4000000000001
$ws.Range("L60").Value = 5534.25
$ws.Range("M60").Value = -264.4000000000001
$ws.Range("N60").Value = -6036.25

$ws.Range("H119").Value = 3999.5
$ws.Range("I119").Value = 1499.25
$ws.Range("K119").Value = 4497.75
$ws.Range("M119").Value = 340.25

$ws.Range("H129").Value = 2108
$ws.Range("J129").Value = 2398.8572
$ws.Range("L129").Value = 7196.571599999999
$ws.Range("N129").Value = -17196.5716

$ws.Range("H130").Value = 586499.7
$ws.Range("I130").Value = 1167999.4
$ws.Range("K130").Value = 3503998.2
$ws.Range("M130").Value = -3498978.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 5009
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 5009
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 5009
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -5235

$ws.Range("H16").Value = 5009
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5009
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5009
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -5509

$ws.Range("H122").Value = 6352.2144
$ws.Range("I122").Value = 5433
$ws.Range("J122").Value = 7271.4287
$ws.Range("K122").Value = 16299
$ws.Range("L122").Value = 21814.2861
$ws.Range("M122").Value = -13849
$ws.Range("N122").Value = -26714.2861

$ws.Range("H132").Value = 5242.5454
$ws.Range("I132").Value = 4050.8
$ws.Range("K132").Value = 12152.4
$ws.Range("M132").Value = -9622.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H122").Value = 50004640
$ws.Range("I122").Value = 250002100
$ws.Range("J122").Value = 5272.4375
$ws.Range("K122").Value = 750006300
$ws.Range("L122").Value = 15817.3125
$ws.Range("M122").Value = -750003850
$ws.Range("N122").Value = -20717.3125

$ws.Range("H132").Value = 4843.1816
$ws.Range("I132").Value = 4140.25
$ws.Range("K132").Value = 12420.75
$ws.Range("M132").Value = -9890.75

$ws.Range("H136").Value = 6630.522
$ws.Range("I136").Value = 5187.6665
$ws.Range("K136").Value = 15562.9995
$ws.Range("M136").Value = -13012.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 10000
$ws.Range("J2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10224

$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1080
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 516.6
$ws.Range("I113").Value = 527.6667
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 1583.0001
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 586.9999
$ws.Range("N113").Value = -5840

$ws.Range("H122").Value = 5999.5
$ws.Range("I122").Value = 4999.25
$ws.Range("K122").Value = 14997.75
$ws.Range("M122").Value = -12547.75

$ws.Range("H126").Value = 2316.9285
$ws.Range("I126").Value = 2386.4167
$ws.Range("K126").Value = 7159.250100000001
$ws.Range("M126").Value = -4689.250100000001

$ws.Range("H132").Value = 2912.2
$ws.Range("J132").Value = 3888.4167
$ws.Range("L132").Value = 11665.2501
$ws.Range("N132").Value = -16725.2501
